# Add files via upload
#
# 1. sharedStrings: "Bangalore" -> "Bangalore-BLR", "Delhi" -> "Delhi-DEL"
#    (these live in Sheet1!B2 and Sheet1!C2)
# 2. Sheet1 becomes the active/selected sheet, with C7 as the selected cell
#    (previously Sheet2 was the active sheet with A10/A2/A10:B17 selected
#    on Sheet1 and B1 selected on Sheet2)
# 3. Sheet2 is no longer the active sheet, but keeps its own B1 selection

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the city values (rewrites the shared-string table entries).
$ws1.Range("B2").Value = "Bangalore-BLR"
$ws1.Range("C2").Value = "Delhi-DEL"

# Make Sheet1 the active sheet and select C7 on it.
[void]$ws1.Activate()
[void]$ws1.Range("C7").Select()
